# Auto-generated edit script: updates leve-profit market-price derived values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3887.0417
$ws.Range("I64").Value = 3753.8333
$ws.Range("J64").Value = 4020.25
$ws.Range("K64").Value = 3753.8333
$ws.Range("L64").Value = 4020.25
$ws.Range("M64").Value = -3505.8333
$ws.Range("N64").Value = -4516.25
$ws.Range("H67").Value = 3887.0417
$ws.Range("I67").Value = 3753.8333
$ws.Range("J67").Value = 4020.25
$ws.Range("K67").Value = 3753.8333
$ws.Range("L67").Value = 4020.25
$ws.Range("M67").Value = -2895.8333
$ws.Range("N67").Value = -5736.25
$ws.Range("H75").Value = 39998
$ws.Range("J75").Value = 39998
$ws.Range("L75").Value = 39998
$ws.Range("N75").Value = -41870
$ws.Range("H78").Value = 39998
$ws.Range("J78").Value = 39998
$ws.Range("L78").Value = 119994
$ws.Range("N78").Value = -129354
$ws.Range("H138").Value = 2292.37
$ws.Range("I138").Value = 789.5925999999999
$ws.Range("J138").Value = 2848.192
$ws.Range("K138").Value = 2368.7778
$ws.Range("L138").Value = 8544.576000000001
$ws.Range("M138").Value = 2771.2222
$ws.Range("N138").Value = -18824.576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1600
$ws.Range("J13").Value = 1600
$ws.Range("L13").Value = 1600
$ws.Range("N13").Value = -1888
$ws.Range("H32").Value = 6231.4
$ws.Range("I32").Value = 5471.066
$ws.Range("J32").Value = 13919.223
$ws.Range("K32").Value = 5471.066
$ws.Range("L32").Value = 13919.223
$ws.Range("M32").Value = -5184.066
$ws.Range("N32").Value = -14493.223
$ws.Range("H61").Value = 3056.0625
$ws.Range("I61").Value = 2598.4546
$ws.Range("J61").Value = 4062.8
$ws.Range("K61").Value = 2598.4546
$ws.Range("L61").Value = 4062.8
$ws.Range("M61").Value = -2386.4546
$ws.Range("N61").Value = -4486.8
$ws.Range("H63").Value = 2899.9167
$ws.Range("J63").Value = 3871.4285
$ws.Range("L63").Value = 3871.4285
$ws.Range("N63").Value = -5243.4285
$ws.Range("H66").Value = 2899.9167
$ws.Range("J66").Value = 3871.4285
$ws.Range("L66").Value = 19357.1425
$ws.Range("N66").Value = -26221.1425
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 1551.3429
$ws.Range("I132").Value = 922.65515
$ws.Range("J132").Value = 4590
$ws.Range("K132").Value = 2767.96545
$ws.Range("L132").Value = 13770
$ws.Range("M132").Value = -237.9654500000001
$ws.Range("N132").Value = -18830
$ws.Range("H136").Value = 3056.0625
$ws.Range("I136").Value = 2598.4546
$ws.Range("J136").Value = 4062.8
$ws.Range("K136").Value = 7795.3638
$ws.Range("L136").Value = 12188.4
$ws.Range("M136").Value = -5245.3638
$ws.Range("N136").Value = -17288.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12465
$ws.Range("I82").Value = 1322.5
$ws.Range("J82").Value = 34750
$ws.Range("K82").Value = 1322.5
$ws.Range("L82").Value = 34750
$ws.Range("M82").Value = -939.5
$ws.Range("N82").Value = -35516
$ws.Range("H85").Value = 12465
$ws.Range("I85").Value = 1322.5
$ws.Range("J85").Value = 34750
$ws.Range("K85").Value = 1322.5
$ws.Range("L85").Value = 34750
$ws.Range("M85").Value = 3.5
$ws.Range("N85").Value = -37402
$ws.Range("H134").Value = 2724.5833
$ws.Range("I134").Value = 2327.4285
$ws.Range("J134").Value = 5504.6665
$ws.Range("K134").Value = 6982.2855
$ws.Range("L134").Value = 16513.9995
$ws.Range("M134").Value = -4447.2855
$ws.Range("N134").Value = -21583.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16537.834
$ws.Range("I28").Value = 3584
$ws.Range("J28").Value = 19128.6
$ws.Range("K28").Value = 3584
$ws.Range("L28").Value = 19128.6
$ws.Range("M28").Value = -3339
$ws.Range("N28").Value = -19618.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 8000050
$ws.Range("I11").Value = 8000050
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 24000150
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -24000010
$ws.Range("N11").ClearContents()
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -2254
$ws.Range("H136").Value = 3074.5
$ws.Range("I136").Value = 1234.1
$ws.Range("J136").Value = 5375
$ws.Range("K136").Value = 3702.3
$ws.Range("L136").Value = 16125
$ws.Range("M136").Value = 1397.7
$ws.Range("N136").Value = -26325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 18855
$ws.Range("J100").Value = 18855
$ws.Range("L100").Value = 18855
$ws.Range("N100").Value = -21019
$ws.Range("H101").Value = 24000
$ws.Range("J101").Value = 24000
$ws.Range("L101").Value = 24000
$ws.Range("N101").Value = -30490
$ws.Range("H102").Value = 1097.3334
$ws.Range("I102").Value = 1194.1538
$ws.Range("J102").Value = 845.6
$ws.Range("K102").Value = 1194.1538
$ws.Range("L102").Value = 845.6
$ws.Range("M102").Value = 427.8462
$ws.Range("N102").Value = -4089.6
$ws.Range("H126").Value = 11967060
$ws.Range("I126").Value = 10001203
$ws.Range("K126").Value = 30003609
$ws.Range("M126").Value = -30001139
$ws.Range("H132").Value = 4486.3335
$ws.Range("I132").Value = 6433.3335
$ws.Range("J132").Value = 2539.3333
$ws.Range("K132").Value = 19300.0005
$ws.Range("L132").Value = 7617.999899999999
$ws.Range("M132").Value = -16770.0005
$ws.Range("N132").Value = -12677.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 694.3158
$ws.Range("I22").Value = 581.25
$ws.Range("J22").Value = 776.5454999999999
$ws.Range("K22").Value = 581.25
$ws.Range("L22").Value = 776.5454999999999
$ws.Range("M22").Value = -286.25
$ws.Range("N22").Value = -1366.5455
$ws.Range("H27").Value = 694.3158
$ws.Range("I27").Value = 581.25
$ws.Range("J27").Value = 776.5454999999999
$ws.Range("K27").Value = 581.25
$ws.Range("L27").Value = 776.5454999999999
$ws.Range("M27").Value = -474.25
$ws.Range("N27").Value = -990.5454999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 48804.2
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 80007
$ws.Range("K9").Value = 2000
$ws.Range("L9").Value = 80007
$ws.Range("M9").Value = -1860
$ws.Range("N9").Value = -80287
$ws.Range("H132").Value = 1859.9683
$ws.Range("I132").Value = 1172.7291
$ws.Range("J132").Value = 4059.1333
$ws.Range("K132").Value = 3518.1873
$ws.Range("L132").Value = 12177.3999
$ws.Range("M132").Value = -988.1873000000001
$ws.Range("N132").Value = -17237.3999
$ws.Range("H136").Value = 4446.923
$ws.Range("I136").Value = 4317.5864
$ws.Range("K136").Value = 12952.7592
$ws.Range("M136").Value = -10402.7592
$ws.Range("H138").Value = 26514.5
$ws.Range("J138").Value = 26514.5
$ws.Range("L138").Value = 26514.5
$ws.Range("N138").Value = -36794.5
